$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.8
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 4.75
$ws.Range("J2").Value = 2.5
$ws.Range("L2").Value = 5.5
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.62
$ws.Range("X2").Value = 7.5
$ws.Range("AB2").Value = 34
$ws.Range("AD2").Value = 6.5
$ws.Range("AN2").Value = 3.6
$ws.Range("AO2").Value = 10
$ws.Range("AX2").Value = 29
$ws.Range("AZ2").Value = 101
